$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.348.12"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.548.97"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'591.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'173.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.88%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "2.548.83"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("D14").Value = "'27.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "3.004.69"
$ws.Range("E15").Value = "  -2.79%  "
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "67.234.70"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "2.569.39"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("D19").Value = "'8.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.43%  "
$ws.Range("D20").Value = "'11.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").Value = "'356.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").Value = "'4.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").Value = "'4.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").Value = "'2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.06%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "'70.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("D27").Value = "'10.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.41%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.668.43"
$ws.Range("E28").Value = "  -2.85%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'0.990"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").Value = "0.0₂01000"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").Value = "'536.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("E32").Value = "  +4.96%  "
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "'157.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("E41").Value = "  -1.82%  "
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("E44").Value = "  +6.35%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'39.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").Value = "'151.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "'0.566"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("D49").Value = "0.0₆0284"
$ws.Range("E49").Value = "  -4.52%  "
$ws.Range("D50").Value = "'3.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("E51").Value = "  +1.35%  "
